$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting C:G to D:H
$ws.Columns("C:C").Insert()

# New header in row 3 for the inserted column
$ws.Range("C3").Value = "Accuracy after attack"
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# New raw-accuracy values for the inserted column (rows 4-28)
$values = @{
    4  = 82.79569892473118
    5  = 5.376344086021505
    6  = 0
    7  = 0
    8  = 93.54838709677421
    9  = 91.39784946236558
    10 = 81.72043010752688
    11 = 49.46236559139785
    12 = 12.90322580645161
    13 = 84.94623655913979
    14 = 59.13978494623656
    15 = 4.301075268817205
    16 = 0
    17 = 0
    18 = 49.46236559139785
    19 = 7.526881720430108
    20 = 7.526881720430108
    21 = 7.526881720430108
    22 = 1.075268817204301
    23 = 7.526881720430108
    24 = 0
    25 = 82.79569892473118
    26 = 6.451612903225806
    27 = 0
    28 = 0
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $values[$row]
    $ws.Cells.Item($row, 2).Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
